$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timing values in rows 5-7 ---
$ws.Range("B5").Value = 0.0008141994476318359
$ws.Range("B6").Value = 0.0005679130554199219
$ws.Range("B7").Value = 0.00701594352722168

# --- Convert tuple-style text to list-style text in A8, A49, A112, A170 ---
$ws.Range("A8").Value = "[[3, 2], [3, 0], [2, 0], [2, 2], [1, 3], [0, 3], [0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [1, 1], [0, 1]]"
$ws.Range("A49").Value = "[[3, 2], [3, 0], [2, 0], [2, 1], [1, 1], [1, 2], [0, 2], [0, 3], [1, 3], [0, 1], [2, 2], [1, 0], [3, 1]]"
$ws.Range("A112").Value = "[[3, 2], [3, 1], [2, 0], [1, 0], [0, 1], [0, 0], [1, 1], [2, 1], [3, 0], [2, 2], [3, 3], [1, 2], [0, 2]]"
$ws.Range("A170").Value = "[[1, 1], [0, 2], [0, 1], [1, 2], [0, 0], [1, 3], [1, 0], [0, 3], [2, 0], [2, 3], [2, 1], [2, 2], [3, 1]]"

# --- Insert a new row before row 215 for the "move_fidelity" entry ---
$ws.Rows.Item(215).Insert()

$ws.Range("A215").Value = "move_fidelity"
$ws.Range("B215").Value = 0.9983505081702928

# --- Update "total time:" row (now shifted to row 219) ---
$ws.Range("B219").Value = 0.03876471519470215
